# Update column header names from "physical quantity" wording to "measurand" wording.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F1").Value = "Measurand"
$ws.Range("I1").Value = "Measurand Level I"
$ws.Range("J1").Value = "Measurand Level II"

# Move the active selection to F2, matching the saved workbook state.
$ws.Range("F2").Select()
